$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins / Losses / Ties), matching the
# existing header styling (bold, centered, bordered) by copying the
# format from the adjacent "Unnamed: 27" header cell.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Every player row (2-45) gets the same season record for this team:
# 96 wins, 66 losses, 0 ties.
$ws.Range("AC2:AC45").Value = 96
$ws.Range("AD2:AD45").Value = 66
$ws.Range("AE2:AE45").Value = 0
